$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 14:58"

# --- Swap country names for rows 42/43 (Suecia <-> Paises Bajos reorder in shared strings) ---
$ws.Range("A42").Value = "Paises Bajos"
$ws.Range("A43").Value = "Suecia"

# --- Update numeric COVID stats across countries ---
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6829951
$ws.Range("C4").Value = 1650
$ws.Range("D4").Value = 4120133
$ws.Range("E4").Value = 2508415
$ws.Range("G4").Value = 55
$ws.Range("H4").Value = 201403

# Row 5 - India
$ws.Range("B5").Value = 5128918
$ws.Range("C5").Value = 13025
$ws.Range("D5").Value = 4029525
$ws.Range("E5").Value = 1016055
$ws.Range("G5").Value = 108
$ws.Range("H5").Value = 83338

# Row 25 - Alemania
$ws.Range("B25").Value = 267182
$ws.Range("C25").Value = 317
$ws.Range("E25").Value = 18631
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 9451

# Row 38 - Belgica
$ws.Range("B38").Value = 97824
$ws.Range("C38").Value = 825
$ws.Range("D38").Value = 87911
$ws.Range("E38").Value = 9338
$ws.Range("G38").Value = 4
$ws.Range("H38").Value = 575

# Row 42 - now Paises Bajos
$ws.Range("B42").Value = 88073
$ws.Range("C42").Value = 1753
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 6266

# Row 43 - now Suecia
$ws.Range("B43").Value = 87885
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 5864

# Row 81 - Libano
$ws.Range("B81").Value = 21393
$ws.Range("C81").Value = 453
$ws.Range("D81").Value = 16918
$ws.Range("E81").Value = 3840
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 635

# Row 85 - Republica de Macedonia
$ws.Range("B85").Value = 16274
$ws.Range("C85").Value = 186
$ws.Range("D85").Value = 13635
$ws.Range("E85").Value = 1964
$ws.Range("G85").Value = 7
$ws.Range("H85").Value = 675

# Row 86 - Madagascar
$ws.Range("B86").Value = 15925
$ws.Range("C86").Value = 54
$ws.Range("D86").Value = 14547
$ws.Range("E86").Value = 1162
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 216

# Row 114 - Eslovaquia
$ws.Range("B114").Value = 5380
$ws.Range("C114").Value = 114
$ws.Range("D114").Value = 2489
$ws.Range("E114").Value = 2831

# Row 194 - Barbados
$ws.Range("B194").Value = 141
$ws.Range("C194").Value = 1
$ws.Range("E194").Value = 5
